$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newString9 = 'U2FsdGVkX1+3dKVevizvc0bfl2VOldCMb0Hk3oGegSXE4NpFNBP8hoavvTsy6dw3IbPP5Fu1Mo6UcVMOuQhRIcSJ6OoQzj6mcdG3yFpQ5YEz4A6rZl3/ldKj+/aZz/pPhXN0BsDNFmwkR9WYewOsmxpyQ0NO+Rts6Qv8pCz+/t9VXo6K5aIqAy3VjjcZ26gidWy2wmFmf8ZHGO758beQ7iqrW8A+WAR5GlLvxE2Bq9TsV43LIzDur/FO8sUJnEmawmoxp0uVo3Ju7HJxipb3rO78j8uE5zsxMi33a0b0gbfMJeOTG3YvF+YlpZoVjg5Siz/TVa6pXAgLvWICmGMkblGtJrwqZ9JRrzia1s8B8DDrsZUZmqocjnwagcQoiguXDN2WDTXy1fGaK/ZgDKjYJOaIr2Yb30Rh5+AxvQDikcY5WFwDm2eQFkyaclk3QDs40KSdPc3V5jAE5mSZ6prEKnH7BkNkChfe+Dhsps3zT3TU3ZLx7V+1PD1N3I/sfx9mM6sC9B+0Rx4ARqBh2e92id0IIl6+GAHjU2bXCQGdmT5d/AZu72JFT28hix9BON3XghhxRLer++rc0xKPMjMb/82AxlkGEVZwPb2eTXiXT6s4gHjfTYu8Bb4V42DrjR7MnBnn7SCOeoTRWFfmSbFIspK+H8BEmDLeEBToVcAol0dI0C4eIPOiwJu4WYO8nCLIwbYKupRZzHdixPYFwCrpyg2Nd6gpiTwBhPxEyjUCRs/EriGP0J+PEffgVA3OuBtbt8tHqzIOeEX8QSEpQp5pseQTy27DE9mAeokjoboXJtwjy8y39hnrnZt71jd059MKjzSvbdfrCWPI8xu0lpuwH/5n33+AbNiNVgbDOLbAkIFgKepPtn9w6pB4Y3CWsHHRfEB+ZmTIUrETSOl6RakN2j1k7g1bgzi553oKQXe0P3haKIWXN8B8iBJq4zmJZNsu3mnypa1HtKQp6v0IXQ7GQgs6sd8x8aAt8oL+Ab3Rh7znNcIt/dWsE6Z9O/PDRiOUYQpAnraYQhpYzt7XAw8eC32dIkYL1g7eikM2ttaaqs1y3J6yMVoDWLxvvAELcREUxqDhZRfECaNjs9jW8RiSw4SrpnO1LSxlXocRLfydt3m4F5/vT1P71gcOEWz/zeYFN0uAMoeStJe504O6+4fA78hHH8H0dMIszQ9YvsokUCMjpKVm2oojbskrZ3fuwEgkwbKTZFPYzKm+WDcisNVlITPVmMD6GzZWTZr4ANNE96BFfRWKLtZiSUxjeaAQZ793M2JokmgxxbwmZDrSQIbKI6q0QfbhRQv68SOmsiNKU3g='
$newString10 = 'U2FsdGVkX19n8cfnv3oO+w3jWYd2WBYQfROoheIgxSp21K06arsG5J4EjKqVUYLAQLz3ravk/4j0kysePbk9JhERUpJwh4xR5VW1grALBMCGAqnLgSzLOt97JVpovoZNu7fU1s0WRSV5uOTOQn10mpAdri6W1mPIKdgdnMnNz/r0ATQ6oilAj9uQ5coyu3ioNrQB2vNacHWLtye/g104eTPfSYVQ18Mq1E/Kwhi70SaNmulyisLQBsZiFVoRo0QZNTXd/TDZS/CDmqws7KiKQIW2eAgToDfPJzfEKkdjfyV99CLwH0Xatmqz2YyCqi/X'

$ws.Range("B2").Value = $newString9
$ws.Range("B4").Value = $newString9
$ws.Range("B6").Value = $newString9

$ws.Range("B3").Value = $newString10
$ws.Range("B5").Value = $newString10
$ws.Range("B7").Value = $newString10
